$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update "Version" value: 0.1.1 -> 0.2.0 ---
$ws1.Range("B3").Value = "0.2.0"

# --- Update "Date" value ---
$ws1.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# --- Insert a new "Jurisdiction" row after "Contact" (row 10), shifting
#     Description/Purpose/Copyright/Immutable down by one row ---
$ws1.Rows.Item(11).Insert()

# Copy formatting (border/fill/alignment) from the row above so the new
# row matches the sheet's existing look instead of Excel's bare default.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = "iso:code:3166:FR"
